$wb = $excel.ActiveWorkbook

# --- Step 1: rename the two existing sheets -------------------------------
$sheetGold = $wb.Worksheets.Item(1)
$sheetGold.Name = "leiden_f gold"

$sheetNotExact = $wb.Worksheets.Item(2)
$sheetNotExact.Name = "leiden_f_goldnotexact"

# --- Step 2: duplicate "leiden_f_goldnotexact" -> "leiden_0_goldnotexact" -
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$null = $sheetNotExact.Copy($null, $last)
$sheet0NotExact = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet0NotExact.Name = "leiden_0_goldnotexact"

# the "leiden_0" run hasn't finished computing every cell yet, so most of
# the pasted-in percentages are still blank
$sheet0NotExact.Range("B2:E4").ClearContents()
$sheet0NotExact.Range("D2").Value = 0.31756756756756699
$sheet0NotExact.Range("E2").Value = 0.43243243243243201
$sheet0NotExact.Range("D4").Value = 0.91216216216216195

# --- Step 3: duplicate "leiden_f gold" -> "leiden_0_gold" -----------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$null = $sheetGold.Copy($null, $last)
$sheet0Gold = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet0Gold.Name = "leiden_0_gold"

$sheet0Gold.Range("B2:E4").ClearContents()

# put the gold tab ahead of the goldnotexact tab
$sheet0Gold.Move($sheet0NotExact)

# the Move() call shuffles worksheet positions, so re-fetch stable handles
# by name before touching either sheet again
$sheet0Gold = $wb.Worksheets.Item("leiden_0_gold")
$sheet0NotExact = $wb.Worksheets.Item("leiden_0_goldnotexact")

# --- Step 4: reset the zoom level on every sheet to 100% normal view ------
$sheetGold.Activate()
$excel.ActiveWindow.Zoom = 100

$sheetNotExact.Activate()
$excel.ActiveWindow.Zoom = 100

$sheet0Gold.Activate()
$null = $sheet0Gold.Range("B18").Select()
$excel.ActiveWindow.Zoom = 100

$sheet0NotExact.Activate()
$null = $sheet0NotExact.Range("C8").Select()
$excel.ActiveWindow.Zoom = 100
